$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MAIN")
$ws2 = $wb.Worksheets.Item("LookUp")

# The old row 8 (the last "educational institution" data row, whose I-column
# carried the tail of the shared formula) is removed; the former row 9 (the
# blank "type" trailer row) shifts up to become the new row 8. Do this on
# both sheets, MAIN and LookUp, since LookUp mirrors MAIN via formulas.

# Drop the conditional-formatting rule that lived on H8 before the rows
# shift (it has no counterpart on the old row 9, so it should disappear
# rather than move).
$ws1.Range("H8").FormatConditions.Delete()

$ws1.Rows("8").Delete()
$ws2.Rows("8").Delete()

# Rename the LookUp sheet to DataLookUp.
$ws2.Name = "DataLookUp"

# Restore the selections left on each sheet.
$ws1.Range("C16").Select()
$ws2.Range("H16").Select()

# DataLookUp becomes the active/selected tab.
$ws2.Activate()
